$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 6552959889.657929
$ws.Range("F2").Value = 2097983131.007991
$ws.Range("L2").Value = 4454976758.649939
$ws.Range("M2").Value = [double]"-2.848260337421415e-08"
$ws.Range("N2").Value = 4454976758.649939
$ws.Range("E17").Value = 4370834007.667166
$ws.Range("F17").Value = 3543707785.959492
$ws.Range("J17").Value = 1490.939999996342
$ws.Range("L17").Value = 333703568.9
$ws.Range("N17").Value = 827126221.7076744
$ws.Range("E18").Value = 903810793.3
$ws.Range("F18").Value = 668126343.45
$ws.Range("L18").Value = 117956879.36
$ws.Range("N18").Value = 235684449.85
$ws.Range("E19").Value = 321818927.0128341
$ws.Range("F19").Value = 215552252.3705085
$ws.Range("J19").Value = 5.060000003658425
$ws.Range("L19").Value = 104986704.01
$ws.Range("N19").Value = 106266674.6423257
$ws.Range("E23").Value = 9723084549.105
$ws.Range("F23").Value = 8264621866.145
$ws.Range("N23").Value = 1458462682.96
$ws.Range("E24").Value = 16258187056.13755
$ws.Range("F24").Value = 13819458996.81999
$ws.Range("N24").Value = 2438728059.317557
$ws.Range("E25").Value = 18226516770.36745
$ws.Range("F25").Value = 15492539253.76501
$ws.Range("N25").Value = 2733977516.602443
$ws.Range("E30").Value = 1610950579.695
$ws.Range("F30").Value = 1369307992.65
$ws.Range("G30").Value = 80547528.98
$ws.Range("I30").Value = 24
$ws.Range("N30").Value = 241642587.0449999
$ws.Range("E33").Value = 15328440827.415
$ws.Range("F33").Value = 12947556688.34
$ws.Range("G33").Value = 750098439.12
$ws.Range("I33").Value = 537
$ws.Range("N33").Value = 2380884139.075
$ws.Range("E48").Value = 683999555.4102353
$ws.Range("F48").Value = 341999777.0789792
$ws.Range("I48").Value = 212.9776999979373
$ws.Range("L48").Value = 103772016.8809517
$ws.Range("M48").Value = 238227761.4503044
$ws.Range("N48").Value = 341999778.3312562
$ws.Range("E50").Value = 337275478.7997647
$ws.Range("F50").Value = 168637739.2510208
$ws.Range("I50").Value = 81.02230000206265
$ws.Range("L50").Value = 50739978.34904829
$ws.Range("M50").Value = 117897761.1996955
$ws.Range("N50").Value = 168637739.5487438
$ws.Range("E55").Value = 8794719092.155001
$ws.Range("F55").Value = 3850585082.325
$ws.Range("L55").Value = 4944134009.83
$ws.Range("M55").Value = [double]"-3.703280526679009e-08"
$ws.Range("N55").Value = 4944134009.83
$ws.Range("E56").Value = 6696996957.425
$ws.Range("F56").Value = 2369283668.3
$ws.Range("L56").Value = 4327713289.125
$ws.Range("M56").Value = [double]"-3.390596248209476e-08"
$ws.Range("N56").Value = 4327713289.125
$ws.Range("E58").Value = 10877315936.18225
$ws.Range("F58").Value = 8802389998.989258
$ws.Range("G58").Value = 1536716993.20816
$ws.Range("K58").Value = 102.2000000000924
$ws.Range("L58").Value = 483413038.78077
$ws.Range("M58").Value = 1352119651.119026
$ws.Range("N58").Value = 2074925937.192997
$ws.Range("E60").Value = 460020251.6399525
$ws.Range("F60").Value = 289893607.9861853
$ws.Range("L60").Value = 170126643.6537672
$ws.Range("M60").Value = [double]"2.938579302732237e-10"
$ws.Range("N60").Value = 170126643.6537672
$ws.Range("E61").Value = 3875736665.699435
$ws.Range("F61").Value = 3014956499.205289
$ws.Range("G61").Value = 686976953.0722671
$ws.Range("K61").Value = 153.7212225696989
$ws.Range("L61").Value = 70176140.27225883
$ws.Range("M61").Value = 238784794.4433374
$ws.Range("N61").Value = 860780166.4941459
$ws.Range("E62").Value = 7123805752.885106
$ws.Range("F62").Value = 4141129948.454005
$ws.Range("L62").Value = 2982675804.431101
$ws.Range("M62").Value = [double]"-1.655924255103602e-08"
$ws.Range("N62").Value = 2982675804.431101
$ws.Range("E63").Value = 5517217092.178309
$ws.Range("F63").Value = 4258005296.585453
$ws.Range("G63").Value = 992790683.2395726
$ws.Range("K63").Value = 134.0787774302087
$ws.Range("L63").Value = 88949454.30697116
$ws.Range("M63").Value = 327748468.9376366
$ws.Range("N63").Value = 1259211795.592857
$ws.Range("E65").Value = 749602718.98
$ws.Range("F65").Value = 374801359.03
$ws.Range("I65").Value = 224
$ws.Range("L65").Value = 74960272.73
$ws.Range("M65").Value = 299841087.22
$ws.Range("N65").Value = 374801359.95
$ws.Range("E66").Value = 1375281515.968
$ws.Range("F66").Value = 753360246.284
$ws.Range("L66").Value = 621921269.684
$ws.Range("N66").Value = 621921269.684
$ws.Range("E67").Value = 5131759767.67
$ws.Range("F67").Value = 2096998391.88
$ws.Range("L67").Value = 3034761375.79
$ws.Range("M67").Value = [double]"-1.90757418749854e-08"
$ws.Range("N67").Value = 3034761375.79
$ws.Range("E69").Value = 722794485.9490119
$ws.Range("F69").Value = 295610293.2928183
$ws.Range("L69").Value = 427184192.6561936
$ws.Range("M69").Value = [double]"-2.253643004218428e-08"
$ws.Range("N69").Value = 427184192.6561936
$ws.Range("E86").Value = 4063673398.916
$ws.Range("F86").Value = 3487426400.184
$ws.Range("G86").Value = 357996281.409
$ws.Range("J86").Value = 453.8
$ws.Range("N86").Value = 576246998.732
$ws.Range("E89").Value = 157505858.804
$ws.Range("F89").Value = 135355697.266
$ws.Range("G89").Value = 7137434.261
$ws.Range("J89").Value = 31.2
$ws.Range("N89").Value = 22150161.538
$ws.Range("E106").Value = 600398633.9557805
$ws.Range("F106").Value = 510338838.4150292
$ws.Range("G106").Value = 65005618.90124533
$ws.Range("N106").Value = 90059795.54075132
$ws.Range("E108").Value = 110526592.5104128
$ws.Range("F108").Value = 55263296.2427064
$ws.Range("I108").Value = 5.499999997283643
$ws.Range("N108").Value = 55263296.26770639
$ws.Range("E109").Value = 4488893788.538008
$ws.Range("F109").Value = 3815559715.468382
$ws.Range("G109").Value = 673334073.0696262
$ws.Range("N109").Value = 673334073.0696262
$ws.Range("E110").Value = 2048969412.158765
$ws.Range("F110").Value = 1662738676.807847
$ws.Range("G110").Value = 386230735.3509181
$ws.Range("N110").Value = 386230735.3509181
$ws.Range("E111").Value = 1628667119.720215
$ws.Range("F111").Value = 1340759489.443177
$ws.Range("G111").Value = 287907630.277038
$ws.Range("N111").Value = 287907630.277038
$ws.Range("E115").Value = 2854299.024405497
$ws.Range("F115").Value = 2426154.136357572
$ws.Range("G115").Value = 428144.8880479246
$ws.Range("N115").Value = 428144.8880479245
$ws.Range("E117").Value = 110526592.6195872
$ws.Range("F117").Value = 55263296.29729361
$ws.Range("I117").Value = 5.500000002716357
$ws.Range("N117").Value = 55263296.32229361
$ws.Range("E118").Value = 102183703.5719921
$ws.Range("F118").Value = 86856147.42161825
$ws.Range("G118").Value = 15327556.15037381
$ws.Range("N118").Value = 15327556.15037381
$ws.Range("E119").Value = 121594046.1791665
$ws.Range("F119").Value = 98673568.24838112
$ws.Range("G119").Value = 22920477.93078539
$ws.Range("N119").Value = 22920477.93078539
$ws.Range("E124").Value = 9226725.199814001
$ws.Range("F124").Value = 7842716.318613252
$ws.Range("G124").Value = 1255577.58070675
$ws.Range("N124").Value = 1384008.88120075
$ws.Range("E128").Value = 58558038.26206814
$ws.Range("F128").Value = 47519847.7537716
$ws.Range("G128").Value = 11038190.50829655
$ws.Range("N128").Value = 11038190.50829655
$ws.Range("E129").Value = 20321571.72978491
$ws.Range("F129").Value = 16729225.75682293
$ws.Range("G129").Value = 3592345.972961977
$ws.Range("N129").Value = 3592345.972961977
$ws.Range("E131").Value = 15435742462.78946
$ws.Range("F131").Value = 12855879852.83287
$ws.Range("G131").Value = 2541830504.804044
$ws.Range("I131").Value = 23.46000000019099
$ws.Range("J131").Value = 29.18999999967338
$ws.Range("K131").Value = 9.340000000895971
$ws.Range("L131").Value = 32756427.76564175
$ws.Range("M131").Value = 2507238.98
$ws.Range("N131").Value = 2579862609.956591
$ws.Range("E135").Value = 5117165920.12
$ws.Range("F135").Value = 4298227197.19
$ws.Range("G135").Value = 469192912.44
$ws.Range("K135").Value = 5
$ws.Range("L135").Value = 349029162.9
$ws.Range("M135").Value = [double]"-1.243279257323593e-09"
$ws.Range("N135").Value = 818938722.9300001
$ws.Range("E137").Value = 8843987119.794807
$ws.Range("F137").Value = 7319666208.751851
$ws.Range("G137").Value = 1153883505.105112
$ws.Range("I137").Value = 166.3999999993564
$ws.Range("J137").Value = 637.0400000065634
$ws.Range("K137").Value = 12.14000000014933
$ws.Range("L137").Value = 51481796.2083945
$ws.Range("M137").Value = [double]"-4.179310053051054e-10"
$ws.Range("N137").Value = 1524320911.042957
$ws.Range("E138").Value = 2246072163.668969
$ws.Range("F138").Value = 1875883495.690736
$ws.Range("G138").Value = 319686794.7782325
$ws.Range("N138").Value = 370188667.9782325
$ws.Range("E139").Value = 448869506.8561808
$ws.Range("F139").Value = 368198976.5401717
$ws.Range("G139").Value = 73474771.54600915
$ws.Range("N139").Value = 80670530.31600915
$ws.Range("E141").Value = 3374709781.360944
$ws.Range("F141").Value = 2696953996.125652
$ws.Range("G141").Value = 639661575.8642932
$ws.Range("J141").Value = 363.7999999972268
$ws.Range("N141").Value = 677755785.2352918
$ws.Range("E145").Value = 14491047188.12485
$ws.Range("F145").Value = 11863017833.97909
$ws.Range("G145").Value = 2443869463.465758
$ws.Range("I145").Value = 95
$ws.Range("J145").Value = 292
$ws.Range("K145").Value = 92
$ws.Range("L145").Value = 85914290.68000001
$ws.Range("M145").Value = 183131867.73
$ws.Range("N145").Value = 2628029354.145758
$ws.Range("E146").Value = 5682167924.56
$ws.Range("F146").Value = 4648445157.04
$ws.Range("G146").Value = 712156812.22
$ws.Range("K146").Value = 552
$ws.Range("L146").Value = 219946095.16
$ws.Range("M146").Value = 467564908.96
$ws.Range("N146").Value = 1033722767.52
$ws.Range("E149").Value = 2535456775.534787
$ws.Range("F149").Value = 2080569923.989627
$ws.Range("G149").Value = 389797940.206551
$ws.Range("I149").Value = 97.14000000045264
$ws.Range("J149").Value = 737.9699999965364
$ws.Range("K149").Value = 187.5199999989547
$ws.Range("L149").Value = 23419720.01596375
$ws.Range("M149").Value = 59814145.24999999
$ws.Range("N149").Value = 454886851.5451602
$ws.Range("E151").Value = 0
$ws.Range("F151").Value = 0
$ws.Range("L151").Value = 0
$ws.Range("M151").Value = 0
$ws.Range("N151").Value = 0
